# Update view-count values in column F across all four sheets
# (data refreshed to a later snapshot; values sourced from commit 456a3b4)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1257
$ws1.Range("F7").Value = 297
$ws1.Range("F8").Value = 1094
$ws1.Range("F10").Value = 6887
$ws1.Range("F14").Value = 7796
$ws1.Range("F16").Value = 47
$ws1.Range("F17").Value = 4949
$ws1.Range("F18").Value = 41
$ws1.Range("F19").Value = 2280
$ws1.Range("F22").Value = 249
$ws1.Range("F23").Value = 365
$ws1.Range("F24").Value = 75
$ws1.Range("F26").Value = 281
$ws1.Range("F27").Value = 230
$ws1.Range("F29").Value = 2010
$ws1.Range("F30").Value = 18
$ws1.Range("F31").Value = 225
$ws1.Range("F32").Value = 67
$ws1.Range("F33").Value = 530
$ws1.Range("F35").Value = 1369
$ws1.Range("F36").Value = 19
$ws1.Range("F37").Value = 2093

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 81
$ws2.Range("F5").Value = 14

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 84

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 84
$ws4.Range("F8").Value = 1257
$ws4.Range("F9").Value = 81
$ws4.Range("F11").Value = 297
$ws4.Range("F12").Value = 1094
$ws4.Range("F14").Value = 6887
$ws4.Range("F18").Value = 7796
$ws4.Range("F20").Value = 47
$ws4.Range("F21").Value = 4949
$ws4.Range("F22").Value = 41
$ws4.Range("F23").Value = 2280
$ws4.Range("F26").Value = 249
$ws4.Range("F27").Value = 365
$ws4.Range("F28").Value = 75
$ws4.Range("F32").Value = 281
$ws4.Range("F33").Value = 230
$ws4.Range("F35").Value = 2010
$ws4.Range("F36").Value = 18
$ws4.Range("F37").Value = 225
$ws4.Range("F38").Value = 67
$ws4.Range("F39").Value = 530
$ws4.Range("F41").Value = 14
$ws4.Range("F42").Value = 1369
$ws4.Range("F43").Value = 19
$ws4.Range("F44").Value = 2093

Write-Output "Updated 50 F-column values across 4 sheets"
